# Refactor synthetic array /3
# The "intervention_type" legend used emoji squares (⬛ noir / 🟥 rouge /
# 🟧 orange / 🟩 vert) as a synthetic color array. This swaps the square
# emoji for book emoji and renames the "noir" (black) label to "bleu"
# (blue), matching the new 📘/📕/📙/📗 legend while leaving rouge/orange/
# vert untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $textA = $cellA.Text

    if ($textA -eq "⬛") {
        $cellA.Value = "📘"
        if ($cellB.Text -eq "noir") {
            $cellB.Value = "bleu"
        }
    } elseif ($textA -eq "🟥") {
        $cellA.Value = "📕"
    } elseif ($textA -eq "🟧") {
        $cellA.Value = "📙"
    } elseif ($textA -eq "🟩") {
        $cellA.Value = "📗"
    }
}
